$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "all": add a new day row (2020-04-28 / serial 43949) above the two
# footnote rows, shifting them down by one and updating the "under
# investigation" count footnote from 35 to 30.
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all")
$wsAll.Rows.Item(21).Insert(-4121)
$wsAll.Range("A21").Value = 43949
$wsAll.Range("B21").Value = 254
$wsAll.Range("C21").Value = 225
$wsAll.Range("D21").Value = 132
$wsAll.Range("E21").Value = 122
$wsAll.Range("F21").Value = 10
$wsAll.Range("G21").Value = 3
$wsAll.Range("H21").Value = 90
$wsAll.Range("B23").Value = "※30件調査中"

# ---------------------------------------------------------------------------
# Sheet "kobe": update the last existing day's follow-up counts, then turn
# the old footer row into a new day row (2020-04-28) and push the footer
# ("入院・入居中") down by one row.
# ---------------------------------------------------------------------------
$wsKobe = $wb.Worksheets.Item("kobe")
$wsKobe.Range("D75").Value = 2
$wsKobe.Range("E75").Value = 253
$wsKobe.Rows.Item(76).Insert(-4121)
$wsKobe.Range("A76").Value = 43949
$wsKobe.Range("B76").Value = 0
$wsKobe.Range("C76").Value = 1772
$wsKobe.Range("D76").Value = 1
$wsKobe.Range("E76").Value = 254
$wsKobe.Range("F76").Value = 127
$wsKobe.Range("G76").Value = 118
$wsKobe.Range("H76").Value = 9
$wsKobe.Range("I76").Value = 3
$wsKobe.Range("J76").Value = 84

# ---------------------------------------------------------------------------
# Sheet "other": turn the old footer row into a new day row (2020-04-28),
# pushing the footer ("※他自治体において...") down by one row and leaving a
# new blank row after it.
# ---------------------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("other")
$wsOther.Rows.Item(51).Insert(-4121)
$wsOther.Range("A51").Value = 43949
$wsOther.Range("B51").Value = 0
$wsOther.Range("C51").Value = 11
$wsOther.Range("D51").Value = 5
$wsOther.Range("E51").Value = 4
$wsOther.Range("F51").Value = 1
$wsOther.Range("G51").Value = 0
$wsOther.Range("H51").Value = 6
